$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 12; $i -le 20; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
    $ws.Cells.Item($i, 2).Value = 0
    $ws.Cells.Item($i, 3).Value = 0
}
